$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$wsBSfVBP = $wb.Worksheets.Item("BSfVBP")

# --- About sheet (sheet1) ---------------------------------------------
# Row 1: title stays the same text, just reshuffled shared-string slot.
$wsAbout.Range("A1").Value = "BSfVBP BAU Subsidy for Vehicle Battery Production"

# Row 3: "Sources:" label stays, but the source itself collapses to "none".
$wsAbout.Range("A3").Value = "Sources:"
$wsAbout.Range("B3").Value = "none"

# Rows 4 & 5 used to hold the publication year + bill title; now blank
# placeholders (formatting retained, content cleared).
$wsAbout.Range("B4").ClearContents()
$wsAbout.Range("B5").ClearContents()

# Row 6 (source URL) is removed entirely.
$wsAbout.Range("B6").ClearContents()

# Row 8: "Notes:" label unchanged.
$wsAbout.Range("A8").Value = "Notes:"

# Rows 9-10: replace the US IRA note with the EU note.
$wsAbout.Range("A9").Value = "In the EU only very specific projects receive funding (e.g. via Projects of Common Interest) but there is no general financial support for any battery production."
$wsAbout.Range("A10").Value = "That is why no financial support is used here. "

# Row 12 (price multiplier used by the BSfVBP formulas) is removed entirely.
$wsAbout.Range("A12:B12").ClearContents()

# --- BSfVBP sheet (sheet2) ---------------------------------------------
# D2:M2 used to compute 45*About!$A$12; now the subsidy no longer applies
# so they become flat literal zeros like the rest of the row.
$wsBSfVBP.Range("D2:M2").Value = 0

# --- Selection / active sheet -------------------------------------------
# The About tab is no longer the active one; BSfVBP becomes active with
# C2 selected, while About keeps A10 selected.
$wsAbout.Range("A10").Select()
$wsBSfVBP.Activate()
$wsBSfVBP.Range("C2").Select()
